# Fix Training Data Issue (#48)
# The "Date" column (BF) held values of the form "4-20-2007-08" which were
# off by one day because of the way the NBA stats site displayed dates.
# Replace them with the corrected ISO-style date text "2008-04-20".
#
# The values must remain plain text (matching the original inlineStr cells),
# so the target number format is forced to Text ("@") before the value is
# written - otherwise Excel will auto-recognize the "2008-04-20" string as a
# real date and silently convert the cell to a date serial number. The cell
# style is then restored to "Normal" so no visible formatting changes are
# introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1

$oldValue = "4-20-2007-08"
$newValue = "2008-04-20"

# Locate the "Date" column from the header row.
$dateCol = 0
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $header = $ws.Cells.Item($firstRow, $col).Value2
    if ($header -eq "Date") {
        $dateCol = $col
        break
    }
}

if ($dateCol -eq 0) {
    # Fallback: workbook is known to use column BF for the Date field.
    $dateCol = 58
}

for ($row = $firstRow + 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}
